$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Years (2010-2020) and the corresponding B:G numeric data extracted from the PDF.
$years = @(2010, 2011, 2012, 2013, 2014, 2015, 2016, 2017, 2018, 2019, 2020)

$data = @(
    @(0, 0, 81, 104, 386, 24),
    @(0, 0, 0, 0, 0, 0),
    @(0, 0, 105, 151, 407, 27),
    @(2, 0, 82, 138, 360, 25),
    @(0, 0, 60, 85, 269, 27),
    @(4, 0, 43, 62, 293, 29),
    @(3, 0, 32, 47, 312, 27),
    @(4, 1, 36, 78, 361, 6),
    @(10, 1, 44, 77, 351, 5),
    @(7, 0, 32, 80, 295, 3),
    @(9, 0, 24, 62, 240, 4)
)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2

    # Column A: the year, stored as a shared-string (text) cell styled like the
    # header row (bold + border + centered), matching the original workbook's
    # convention without introducing any new cell-format (xf) entries.
    #
    # Route the number through TEXT() so the resulting literal is a real text
    # value, then "paste values" to bake it in as plain text (t="s"), and
    # finally "paste formats" from a header cell to reuse its existing style
    # index instead of minting a brand-new one (e.g. via NumberFormat, which
    # would otherwise create an unused extra xf entry).
    $ws.Cells.Item($row, 1).Formula = '=TEXT(' + $years[$i] + ',"0")'
    $ws.Cells.Item($row, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4163)

    $ws.Cells.Item(1, 2).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    # Columns B:G: plain numeric values, unstyled (matches the original cells).
    for ($j = 0; $j -lt 6; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}

$excel.CutCopyMode = $false
